# Fix "gurantee" typo / account field-name typos in the diagram shapes:
#   start_time -> startTime
#   end_time   -> endTime
# (and, best-effort, the notes-master auto date field 2023/3/10 -> 2023/4/27)

function Replace-InTextRange {
    param($tr, $oldStr, $newStr)

    $text = $tr.Text
    $searchFrom = 0
    while ($true) {
        $idx = $text.IndexOf($oldStr, $searchFrom)
        if ($idx -lt 0) {
            break
        }
        $sub = $tr.Characters($idx + 1, $oldStr.Length)
        $sub.Text = $newStr
        $text = $tr.Text
        $searchFrom = $idx + $newStr.Length
    }
}

function Fix-ShapeTree {
    param($shapes)

    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)

        if ($shp.HasTextFrame) {
            if ($shp.TextFrame.HasText) {
                $tr = $shp.TextFrame.TextRange
                $txt = $tr.Text
                if ($txt.IndexOf("start_time") -ge 0) {
                    Replace-InTextRange $tr "start_time" "startTime"
                }
                $txt2 = $tr.Text
                if ($txt2.IndexOf("end_time") -ge 0) {
                    Replace-InTextRange $tr "end_time" "endTime"
                }
            }
        }

        # ppGroup = 6 : recurse into grouped shapes so nested diagram text boxes
        # (the "start_time"/"end_time" labels live several group-levels deep) get fixed too.
        if ($shp.Type -eq 6) {
            Fix-ShapeTree $shp.GroupItems
        }
    }
}

$p = $ppt.ActivePresentation

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $s = $p.Slides.Item($si)
    Fix-ShapeTree $s.Shapes
}

# Best-effort: the notes master's auto-updating date field text also shifted
# (2023/3/10 -> 2023/4/27) in the source commit. Attempt the standard
# PowerPoint COM route for it; harmless no-op if the host can't persist it.
try {
    $nm = $p.NotesMaster
    $dt = $nm.HeadersFooters.DateAndTime
    $dt.Value = "2023/4/27"
} catch {
}
